$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts existing B..N to C..O)
$ws.Columns("B").Insert()

# New column B width (matches the author's narrower "tienda" column)
$ws.Columns("B").ColumnWidth = 11.6666666666

# Header for the newly inserted column
$ws.Range("B1").Value = "tienda"

# Fill the ten data rows with tienda1..tienda10
for ($i = 1; $i -le 10; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = "tienda$i"
}

# Restore a plain view/selection state (matches the saved file's UI state)
$ws.Range("C16:C18").Select()
